$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column H (rows 2-11) with the new 2023 value of 171.5
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = 171.5
}

# H12 becomes a SUM formula totaling H2:H11
$ws.Range("H12").Formula = "=SUM(H2:H11)"

# Update the selected cell in the sheet view
$ws.Range("I15").Select()
